# Update two test-data values on the "Testdata" sheet (row 11/12 of the
# "Add Course" test case, inside Table2 which spans A10:E16):
#   D12 ("Category"): "Testing"             -> "Automation"
#   D11 ("Name"):     "Selenium Testing987" -> "selenium Testing987"
#
# Note: D12 is written before D11 so that new shared-string entries are
# appended to the shared strings table in the same order as the source
# workbook (i.e. "Automation" before "selenium Testing987").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Testdata")

$ws.Range("D12").Value = "Automation"
$ws.Range("D11").Value = "selenium Testing987"
